$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix ItemSlot: ChainSaw damage reduced (row 14, G) ---
$ws.Range("G14").Value2 = 25

# --- Add Weapon: row 16 becomes "ForestMourn" (a TwohandSword) ---
$ws.Range("C16").Value2 = "ForestMourn"
$ws.Range("E16").Value2 = 2000
$ws.Range("F16").Value2 = "TwohandSword"
$ws.Range("G16").Value2 = 100
$ws.Range("H16").Value2 = 1.5

# --- Add Weapon: row 17 becomes "Mace" ---
$ws.Range("C17").Value2 = "Mace"
$ws.Range("F17").Value2 = "Mace"
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = 0.4

# --- Remove placeholder data from rows 18-20 (now empty weapon slots) ---
$ws.Range("D18:F18").ClearContents()
$ws.Range("D19:F19").ClearContents()
$ws.Range("D20:F20").ClearContents()

# --- Fix UI: update selection to G17 ---
$ws.Activate()
$ws.Range("G17").Select()

$wb.Save()
